# Updates the 广州-漫展信息 workbook to the "output generated at 456a3b4" snapshot.
#
# Changes:
#  1. Bump a handful of "想去人数" (F column) counters across the
#     展览 / 演出 / 本地生活 / 全部类型 sheets.
#  2. Insert a newly-scraped event ("广州·卡农·世界经典音乐之旅交响音乐会",
#     2024-10-27) into the 演出 sheet (as row 17) and into the 全部类型
#     sheet (as row 41), pushing the later rows down by one.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, [string]$text)
    # Plain string assignment lets Excel "helpfully" reinterpret things that
    # look like dates (e.g. "2024-10-27") as date serials. Force the cell to
    # Text first so the literal string is preserved, then drop the style
    # back to Normal so we don't leave a stray number format behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Bump-F {
    param($ws, [int]$row, [double]$newValue)
    $ws.Cells.Item($row, 6).Value = $newValue
}

function Insert-KanonRow {
    param($ws, [int]$row, [int]$oldLastRow)

    # Push existing rows down, inheriting the row-above's formatting.
    $ws.Rows.Item($row).Insert()

    # Column A keeps the bordered/bold index style used throughout the
    # sheet; copy that formatting explicitly (Insert() alone leaves a
    # slightly different auto style) onto the freshly inserted row.
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    Set-TextValue $ws.Cells.Item($row, 2) "2024-10-27"
    $ws.Cells.Item($row, 3).Value = "广州·卡农·世界经典音乐之旅交响音乐会"
    $ws.Cells.Item($row, 4).Value = "东风中路299号 广州中山纪念堂"
    $ws.Cells.Item($row, 5).Value = "2024.10.27 19:30-10.27 21:00"
    $ws.Cells.Item($row, 6).Value = 0
    $ws.Cells.Item($row, 7).Value = 75
    $ws.Cells.Item($row, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91040"
    $ws.Cells.Item($row, 9).Value = "//i1.hdslb.com/bfs/openplatform/202408/WEqD8aj31724134831558.jpeg"

    # The source site renumbers column A (0-based running index) on every
    # scrape; it is NOT simply dragged down with the rest of the row, so
    # every row from the insertion point through the old last row keeps
    # counting up by one (the brand new last row continues the sequence).
    for ($r = $row; $r -le ($oldLastRow + 1); $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

# ---------------------------------------------------------------------
# Sheet: 展览 (exhibitions)
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
Bump-F $wsExpo 3  26479
Bump-F $wsExpo 4  584
Bump-F $wsExpo 5  251
Bump-F $wsExpo 6  599
Bump-F $wsExpo 8  540
Bump-F $wsExpo 10 355
Bump-F $wsExpo 11 227
Bump-F $wsExpo 12 188
Bump-F $wsExpo 14 296
Bump-F $wsExpo 15 55
Bump-F $wsExpo 16 392
Bump-F $wsExpo 17 56
Bump-F $wsExpo 18 1518
Bump-F $wsExpo 19 191

# ---------------------------------------------------------------------
# Sheet: 演出 (performances) — F bumps + new row 17
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
Bump-F $wsShow 6  193
Bump-F $wsShow 10 436
Bump-F $wsShow 12 7
Bump-F $wsShow 13 11
Bump-F $wsShow 15 56

Insert-KanonRow $wsShow 17 21

# ---------------------------------------------------------------------
# Sheet: 本地生活 (local life)
# ---------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
Bump-F $wsLocal 2 5026
Bump-F $wsLocal 3 219

# ---------------------------------------------------------------------
# Sheet: 全部类型 (all types) — F bumps + new row 41
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
Bump-F $wsAll 3  5026
Bump-F $wsAll 4  219
Bump-F $wsAll 5  26479
Bump-F $wsAll 6  584
Bump-F $wsAll 8  251
Bump-F $wsAll 10 599
Bump-F $wsAll 14 193
Bump-F $wsAll 15 193
Bump-F $wsAll 19 436
Bump-F $wsAll 20 540
Bump-F $wsAll 23 355
Bump-F $wsAll 24 227
Bump-F $wsAll 25 188
Bump-F $wsAll 27 7
Bump-F $wsAll 28 296
Bump-F $wsAll 29 55
Bump-F $wsAll 30 11
Bump-F $wsAll 32 392
Bump-F $wsAll 33 56
Bump-F $wsAll 34 56
Bump-F $wsAll 35 1518
Bump-F $wsAll 36 191

Insert-KanonRow $wsAll 41 46
